$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh (GitHub Actions data pull) — update price/volume(1h) cells,
# and swap the ARBITRUM / TrustWalletToken rows (41-42) which changed rank order.

# Row 2: D2, E2
$ws.Range("D2").Value = '''25.852.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.41%  '

# Row 3: D3, E3
$ws.Range("D3").Value = '''1.586.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.20%  '

# Row 4: E4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5: D5, E5
$ws.Range("D5").Value = '''209.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.28%  '

# Row 6: E6
$ws.Range("E6").Value = '  +0.05%  '

# Row 7: D7, E7
$ws.Range("D7").Value = '''0.479'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.66%  '

# Row 8: E8
$ws.Range("E8").Value = '  -0.85%  '

# Row 9: E9
$ws.Range("E9").Value = '  -0.32%  '

# Row 10: D10, E10
$ws.Range("D10").Value = '''18.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '

# Row 11: D11, E11
$ws.Range("D11").Value = '''0.0791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '

# Row 12: D12, E12
$ws.Range("D12").Value = '''1.805.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.22%  '

# Row 13: D13, E13
$ws.Range("D13").Value = '''1.587.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.00%  '

# Row 14: E14
$ws.Range("E14").Value = '  -2.81%  '

# Row 15: E15
$ws.Range("E15").Value = '  -2.93%  '

# Row 16: D16, E16
$ws.Range("D16").Value = '''25.834.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '

# Row 18: D18, E18
$ws.Range("D18").Value = '''59.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.19%  '

# Row 19: E19
$ws.Range("E19").Value = '  -0.01%  '

# Row 20: D20, E20
$ws.Range("D20").Value = '''191.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.16%  '

# Row 21: E21
$ws.Range("E21").Value = '  -1.77%  '

# Row 22: E22
$ws.Range("E22").Value = '  -1.84%  '

# Row 23: D23, E23
$ws.Range("D23").Value = '''5.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.49%  '

# Row 24: E24
$ws.Range("E24").Value = '  -0.89%  '

# Row 25: D25, E25
$ws.Range("D25").Value = '''142.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.19%  '

# Row 26: E26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27: E27
$ws.Range("E27").Value = '  -0.66%  '

# Row 28: D28, E28
$ws.Range("D28").Value = '''15.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.94%  '

# Row 29: D29, E29
$ws.Range("D29").Value = '''6.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.91%  '

# Row 30: E30
$ws.Range("E30").Value = '  -5.43%  '

# Row 31: D31, E31
$ws.Range("D31").Value = '''0.0470'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.00%  '

# Row 32: E32
$ws.Range("E32").Value = '  -0.33%  '

# Row 33: E33
$ws.Range("E33").Value = '  -2.31%  '

# Row 34: D34, E34
$ws.Range("D34").Value = '''1.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.18%  '

# Row 35: D35, E35
$ws.Range("D35").Value = '''2.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.39%  '

# Row 36: D36, E36
$ws.Range("D36").Value = '''1.101.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.29%  '

# Row 37: E37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38: E38
$ws.Range("E38").Value = '  -2.20%  '

# Row 39: E39
$ws.Range("E39").Value = '  -2.53%  '

# Row 40: E40
$ws.Range("E40").Value = '  -2.12%  '

# Row 41: B41, C41, D41, E41
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.777'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.12%  '

# Row 42: B42, C42, D42, E42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.818'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.83%  '

# Row 43: D43, E43
$ws.Range("D43").Value = '''5.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.59%  '

# Row 44: D44, E44
$ws.Range("D44").Value = '''93.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.93%  '

# Row 45: D45, E45
$ws.Range("D45").Value = '''1.719.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.15%  '

# Row 46: D46, E46
$ws.Range("D46").Value = '''0.0₆0107'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.52%  '

# Row 47: D47
$ws.Range("D47").Value = '''1.51'
$ws.Range("D47").Style = "Normal"

# Row 48: D48, E48
$ws.Range("D48").Value = '''53.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.67%  '

# Row 49: E49
$ws.Range("E49").Value = '  -1.58%  '

# Row 50: D50, E50
$ws.Range("D50").Value = '''0.408'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.59%  '

# Row 51: E51
$ws.Range("E51").Value = '  -0.08%  '
